# Applies the updates described by the commit:
# "Added UserList test cases of ABP and modified CMCDeferredCorp in IWP"
#
# Concretely (per the OOXML diff for MultibillCCData.xlsx): a batch of
# Katalon test runs were recorded. Several sheets got a new run timestamp
# written into their "Date" column (B), and a handful of those runs also
# flipped their "Result" column (A) from Pass to Fail.

$wb = $excel.ActiveWorkbook

# Each row: SheetName, DateCell, DateValue, ResultCell (or ""), ResultValue (or "")
$changes = @(
    @("VerifyAmountTextBoxEditable", "B2", "Fri Jul 11 23:48:47 IST 2025", "A2", "Fail"),
    @("VerifyLookup1Search",         "B2", "Fri Jul 11 23:50:00 IST 2025", "",   ""),
    @("VerifyStreetAddressSearch",   "B2", "Fri Jul 11 23:50:40 IST 2025", "",   ""),
    @("VerifyUDF3Saerch",            "B2", "Fri Jul 11 23:51:25 IST 2025", "",   ""),
    @("VerifyPaymentEntryPageCorp",  "B2", "Fri Jul 11 23:55:48 IST 2025", "A2", "Fail"),
    @("VerifyPaymentEntryPageCC",    "B2", "Fri Jul 11 23:56:08 IST 2025", "A2", "Fail"),
    @("VerifyPaymentEntryPagePC",    "B2", "Fri Jul 11 23:56:33 IST 2025", "A2", "Fail"),
    @("VerifyRemoveCartContent",     "B2", "Fri Jul 11 23:56:55 IST 2025", "A2", "Fail"),
    @("VerifySearchResult",          "B2", "Fri Jul 11 23:57:47 IST 2025", "",   ""),
    @("VerifyStaticTextOnViewCart",  "B2", "Fri Jul 11 23:58:32 IST 2025", "",   ""),
    @("VerifyStaticTextOnViewCart",  "B3", "Fri Jul 11 23:59:13 IST 2025", "",   ""),
    @("VerifyStaticTextOnSearch",    "B2", "Fri Jul 11 23:59:53 IST 2025", "",   ""),
    @("Verify2Pages",                "B2", "Sat Jul 12 00:00:32 IST 2025", "A2", "Fail")
)

foreach ($c in $changes) {
    $sheetName   = $c[0]
    $dateCell    = $c[1]
    $dateValue   = $c[2]
    $resultCell  = $c[3]
    $resultValue = $c[4]

    $ws = $wb.Worksheets.Item($sheetName)

    if ($resultCell -ne "") {
        $ws.Range($resultCell).Value = $resultValue
    }

    $ws.Range($dateCell).Value = $dateValue
}

$wb.Save()
